$d = $word.ActiveDocument

# Change 1: update the "recherche d'alternance" sentence, preserving the
# untouched leading portion as its own run.
$d.Content.Find.Execute(
    "à raison de 3 jours par semaine en entreprise et deux en formation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "dont les modalités de présence en entreprise et en formation sont flexibles, adaptables aux besoins de l'employeur.",
    2
)
